$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column header in H1, reusing the same header style as the
# other header cells (e.g. G1: bold, bordered, centered) by copying G1's
# formatting onto H1, then overwriting the copied text with the new label.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 (plain, unstyled like the other
# numeric data cells in row 2).
$ws.Range("H2").Value = 1
